$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of trade data (row 4) matching the style/format of existing rows
$ws.Range("A4").Value = 10013
$ws.Range("B4").Value = 9998
$ws.Range("C4").Value = 80.45
$ws.Range("D4").Value = 80.569999999999993
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.15

# Copy the date-formatted cell's format onto G4, then set its value
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("G4").Value = 42608.624085648145

$ws.Range("H4").Value = $true
